$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'331.97"
$ws.Range("E2").Value = "'1.41%"
$ws.Range("D3").Value = "'44.65"
$ws.Range("E3").Value = "'1.17%"
$ws.Range("D4").Value = "'5.564"
$ws.Range("E4").Value = "'-0.11%"
$ws.Range("D5").Value = "'0.08284"
$ws.Range("E5").Value = "'2.74%"
$ws.Range("D6").Value = "'2.041"
$ws.Range("E6").Value = "'3.57%"
$ws.Range("D7").Value = "'0.9767"
$ws.Range("E7").Value = "'3.32%"
$ws.Range("D8").Value = "'0.1122"
$ws.Range("E8").Value = "'-3.93%"
$ws.Range("D9").Value = "'0.1914"
$ws.Range("E9").Value = "'2.47%"
$ws.Range("E10").Value = "'-12.84%"
$ws.Range("D11").Value = "'0.1006"
$ws.Range("E11").Value = "'2.14%"
$ws.Range("D12").Value = "'0.04682"
$ws.Range("E12").Value = "'-1.35%"
$ws.Range("E13").Value = "'-0.71%"
$ws.Range("D14").Value = "'0.001271"
$ws.Range("E14").Value = "'-1.86%"
$ws.Range("D15").Value = "'0.04103"
$ws.Range("E15").Value = "'-2.63%"
$ws.Range("D16").Value = "'0.006090"
$ws.Range("E16").Value = "'3.47%"
$ws.Range("D17").Value = "'3.362"
$ws.Range("E17").Value = "'-0.28%"
$ws.Range("D18").Value = "'4.438"
$ws.Range("E18").Value = "'2.73%"
$ws.Range("D20").Value = "'0.3353"
$ws.Range("E20").Value = "'-3.51%"
$ws.Range("D21").Value = "'0.1386"
$ws.Range("E21").Value = "'-2.49%"
$ws.Range("D23").Value = "'0.001302"
$ws.Range("E23").Value = "'3.76%"
$ws.Range("D24").Value = "'0.004393"
$ws.Range("E24").Value = "'2.22%"
$ws.Range("D25").Value = "'0.0001280"
$ws.Range("E25").Value = "'7.34%"
$ws.Range("D26").Value = "'0.0003740"
$ws.Range("E26").Value = "'-0.31%"
$ws.Range("D38").Value = "'0.02796"
$ws.Range("E38").Value = "'7.95%"
$ws.Range("D39").Value = "'0.05745"
$ws.Range("E39").Value = "'4.07%"
$ws.Range("D40").Value = "'0.007630"
$ws.Range("E40").Value = "'0.81%"
$ws.Range("D41").Value = "'0.1422"
$ws.Range("E41").Value = "'1.60%"
$ws.Range("D42").Value = "'0.007567"
$ws.Range("E42").Value = "'1.05%"
$ws.Range("D43").Value = "'0.001974"
$ws.Range("E43").Value = "'-2.29%"
$ws.Range("D44").Value = "'0.008315"
$ws.Range("E44").Value = "'-0.57%"
$ws.Range("D45").Value = "'0.00007041"
$ws.Range("E45").Value = "'-0.78%"
$ws.Range("E46").Value = "'-0.20%"
$ws.Range("D47").Value = "'0.0005807"
$ws.Range("E47").Value = "'-0.08%"
$ws.Range("D48").Value = "'0.003590"
$ws.Range("E48").Value = "'-25.79%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.20%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.20%"
